# Convert the plain-text "Youtube: https://youtu.be/_d68O7wgSDI" line into a HYPERLINK
# field (fldChar begin/instrText/separate/end), matching the Github/Circle CI links above,
# and pointing at the new video URL https://youtu.be/jxgBqiDfKnE.

$d = $word.ActiveDocument

# Find the paragraph holding the old plain-text Youtube link.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Youtube:*") {
        $target = $p
    }
}

$searchRng = $d.Range($target.Range.Start, $target.Range.End)
$find = $searchRng.Find
$find.Execute("Youtube: https://youtu.be/_d68O7wgSDI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Build a fresh Range from the found Start/End (Range objects returned directly off a
# Find don't replace correctly via InsertXML).
$replaceRng = $d.Range($searchRng.Start, $searchRng.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Youtube: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "https://youtu.be/jxgBqiDfKnE" </w:instrText></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="4"/><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>https://youtu.be</w:t></w:r><w:r><w:rPr><w:rStyle w:val="4"/><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>/jxgBqiDfKnE</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$replaceRng.InsertXML($xml)

# InsertXML re-anchors the pre-existing "_GoBack" bookmark across the whole replaced
# span instead of keeping its original (collapsed) position; drop it and recreate it
# precisely between the "https://youtu.be" run and the "/jxgBqiDfKnE" run.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$para = $target
$paraRng = $d.Range($para.Range.Start, $para.Range.End)
$find2 = $paraRng.Find
$find2.Execute("https://youtu.be", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$bmRng = $d.Range($paraRng.End, $paraRng.End)
$d.Bookmarks.Add("_GoBack", $bmRng)
